$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.876175666666667
$ws.Range("H2").Value = 5.628527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.063418
$ws.Range("N2").Value = 3.190254
$ws.Range("O2").Value = 0.0908628824358815
$ws.Range("P2").Value = 0.1028742253445915
$ws.Range("Q2").Value = 1.995158975095333
$ws.Range("R2").Value = 17.956430775858
$ws.Range("S2").Value = 0.0908628824358815
$ws.Range("T2").Value = 0.1028742253445915

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.876175666666667
$ws.Range("H3").Value = 5.628527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.853215333333333
$ws.Range("N3").Value = 5.559646
$ws.Range("O3").Value = 0.1583464704951765
$ws.Range("P3").Value = 0.1792786014656379
$ws.Range("Q3").Value = 3.476957513493555
$ws.Range("R3").Value = 31.292617621442
$ws.Range("S3").Value = 0.1583464704951765
$ws.Range("T3").Value = 0.1792786014656379

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.876175666666667
$ws.Range("H4").Value = 5.628527
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.214728
$ws.Range("N4").Value = 6.644184
$ws.Range("O4").Value = 0.1892356250237018
$ws.Range("P4").Value = 0.2142510540060226
$ws.Range("Q4").Value = 4.155218781885333
$ws.Range("R4").Value = 37.396969036968
$ws.Range("S4").Value = 0.1892356250237018
$ws.Range("T4").Value = 0.2142510540060226

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.876175666666667
$ws.Range("H5").Value = 5.628527
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.472753
$ws.Range("N5").Value = 7.418259000000001
$ws.Range("O5").Value = 0.2112823604001186
$ws.Range("P5").Value = 0.239212190637656
$ws.Range("Q5").Value = 4.639319008277001
$ws.Range("R5").Value = 41.75387107449301
$ws.Range("S5").Value = 0.2112823604001186
$ws.Range("T5").Value = 0.239212190637656

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.876175666666667
$ws.Range("H6").Value = 5.628527
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 4.0994325
$ws.Range("N6").Value = 8.198865
$ws.Range("O6").Value = 0.3502726616451215
$ws.Range("P6").Value = 0.2643839285460921
$ws.Range("Q6").Value = 7.6912555036425
$ws.Range("R6").Value = 46.147533021855
$ws.Range("S6").Value = 0.3502726616451215
$ws.Range("T6").Value = 0.2643839285460921
